$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (B, C): simple string assignment is safe (non-numeric content)
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

# Price column (D): force text format to preserve exact string (avoid numeric auto-conversion),
# then reset style back to Normal so no stray style index is introduced.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '42.678.53'
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.242.14'
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '113.00'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '294.41'
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.628'
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.604'
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '43.91'
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0923'
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '54.47'
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '8.84'
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '1.07'
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '15.02'
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.577.83'
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '2.241.35'
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '42.621.89'
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.23'
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.0000106'
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '74.48'
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '249.27'
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.38'
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '8.96'
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '11.48'
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.21'
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '175.36'
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '37.32'
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '21.76'
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.0885'
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '5.70'
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '5.02'
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '4.23'
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0374'
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.41'
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '71.36'
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.230'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '12.40'
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.31'
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '5.46'
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '104.36'
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '8.55'
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0978'
$c.Style = "Normal"

# Volume(1h) column (E): percentage text with padding spaces, safe as plain string assignment
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("E6").Value = '  +5.18%  '
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("E8").Value = '  -0.36%  '
$ws.Range("E9").Value = '  -0.33%  '
$ws.Range("E10").Value = '  -4.60%  '
$ws.Range("E11").Value = '  -0.55%  '
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("E14").Value = '  +23.05%  '
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("E16").Value = '  -1.63%  '
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("E19").Value = '  -0.76%  '
$ws.Range("E20").Value = '  +7.27%  '
$ws.Range("E21").Value = '  -1.14%  '
$ws.Range("E22").Value = '  +3.47%  '
$ws.Range("E23").Value = '  +9.26%  '
$ws.Range("E24").Value = '  +7.77%  '
$ws.Range("E25").Value = '  +2.25%  '
$ws.Range("E26").Value = '  -2.87%  '
$ws.Range("E27").Value = '  -0.81%  '
$ws.Range("E28").Value = '  -5.43%  '
$ws.Range("E29").Value = '  -1.31%  '
$ws.Range("E30").Value = '  +1.08%  '
$ws.Range("E31").Value = '  -7.27%  '
$ws.Range("E32").Value = '  +3.30%  '
$ws.Range("E33").Value = '  -4.29%  '
$ws.Range("E34").Value = '  -1.85%  '
$ws.Range("E35").Value = '  +2.58%  '
$ws.Range("E36").Value = '  +7.97%  '
$ws.Range("E37").Value = '  -1.35%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("E39").Value = '  +1.01%  '
$ws.Range("E40").Value = '  -1.66%  '
$ws.Range("E41").Value = '  -5.72%  '
$ws.Range("E42").Value = '  +0.76%  '
$ws.Range("E43").Value = '  -0.38%  '
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("E45").Value = '  -5.60%  '
$ws.Range("E46").Value = '  -1.67%  '
$ws.Range("E47").Value = '  -3.21%  '
$ws.Range("E48").Value = '  +1.44%  '
$ws.Range("E49").Value = '  +4.18%  '
$ws.Range("E50").Value = '  +1.93%  '
$ws.Range("E51").Value = '  -0.81%  '
